$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data row (row 2) of its values, except the two styled
# placeholder cells I2/J2 which should keep their style but lose content.
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# Update the selected range shown in the worksheet view.
$ws.Range("A2:AB12").Select() | Out-Null
